{"js": "// The document contains several \"<id>...</id>\" markers, each originally split\n// across three runs: one run for the literal \"<id>\" text, a middle run for the\n// identifier value (e.g. \"p169v_1\"), and a trailing run for the literal\n// \"</id>\" text. For the newly downloaded \"tc\"/\"tcn\" pages (p169v_1 and\n// p169v_2) the three runs should be merged into a single run (keeping the\n// formatting of the opening \"<id>\" run) whose text is the full\n// \"<id>p169v_1</id>\" / \"<id>p169v_2</id>\" string.\nconst ids = [\"p169v_1\", \"p169v_2\"];\n\nfor (const id of ids) {\n  const target = \"<id>\" + id + \"</id>\";\n  const results = context.document.body.search(target, {\n    matchCase: true,\n    matchWildcards: false,\n  });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    // Replacing the found range (which spans the original 3 runs) with the\n    // same text collapses it into one run that carries the formatting of\n    // the first of the originally-matched runs.\n    results.items[i].insertText(target, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# The document contains several \"<id>...</id>\" markers, each originally split\n# across three runs: one run for the literal \"<id>\" text, a middle run for the\n# identifier value (e.g. \"p169v_1\"), and a trailing run for the literal\n# \"</id>\" text. For the newly downloaded \"tc\"/\"tcn\" pages (p169v_1 and\n# p169v_2) the three runs should be merged into a single run (keeping the\n# formatting of the opening \"<id>\" run) whose text is the full\n# \"<id>p169v_1</id>\" / \"<id>p169v_2</id>\" string.\n$d = $word.ActiveDocument\n$openTag = \"<id>\"\n$ids = @(\"p169v_1\", \"p169v_2\")\n\nforeach ($id in $ids) {\n    $fullTarget = $openTag + $id + \"</id>\"\n\n    # Locate the \"<id>...</id>\" span; it currently spans the 3 original runs.\n    $matchRange = $d.Content\n    $found = $matchRange.Find.Execute($fullTarget, $false, $false, $false, $false, $false, $true, 1, $false, \"\", 0)\n\n    if ($found) {\n        $openStart = $matchRange.Start\n        $openEnd = $openStart + $openTag.Length\n        $fullEnd = $matchRange.End\n\n        # Remove the old middle run (\"p169v_1\") and closing run (\"</id>\"),\n        # then retype their combined text right after what remains of the\n        # opening \"<id>\" run, so the result keeps that first run's\n        # formatting/identity.\n        $tailRange = $d.Range($openEnd, $fullEnd)\n        $tailRange.Delete()\n\n        $openRunRange = $d.Range($openStart, $openEnd)\n        $openRunRange.InsertAfter($id + \"</id>\")\n    }\n}\n"}
